$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1050
$ws1.Range("F3").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1050
$ws4.Range("F3").Value = 20
